$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 553. Existing rows 553:571
# shift down to become 555:573 (Excel carries formatting down automatically,
# same as a manual row insert).
$ws.Rows("553:554").Insert()

# Row 553 (new) - "Primera" record dated 45075
$ws.Cells.Item(553, 1).Value = 5
$ws.Cells.Item(553, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(553, 3).Value = "Maule"
$ws.Cells.Item(553, 4).Value = 45075
$ws.Cells.Item(553, 5).Value = 7
$ws.Cells.Item(553, 6).Value = 100114014
$ws.Cells.Item(553, 7).Value = "Betarraga"
$ws.Cells.Item(553, 8).Value = "Sin especificar"
$ws.Cells.Item(553, 9).Value = "Primera"
$ws.Cells.Item(553, 10).Value = 3000
$ws.Cells.Item(553, 11).Value = 600
$ws.Cells.Item(553, 12).Value = 600
$ws.Cells.Item(553, 13).Value = 600
$ws.Cells.Item(553, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(553, 15).Value = "Región del Maule"
$ws.Cells.Item(553, 16).Value = 120
$ws.Cells.Item(553, 17).Value = 5
$ws.Cells.Item(553, 18).Value = "Hortaliza"

# Row 554 (new) - "Segunda" record dated 45075
$ws.Cells.Item(554, 1).Value = 5
$ws.Cells.Item(554, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(554, 3).Value = "Maule"
$ws.Cells.Item(554, 4).Value = 45075
$ws.Cells.Item(554, 5).Value = 7
$ws.Cells.Item(554, 6).Value = 100114014
$ws.Cells.Item(554, 7).Value = "Betarraga"
$ws.Cells.Item(554, 8).Value = "Sin especificar"
$ws.Cells.Item(554, 9).Value = "Segunda"
$ws.Cells.Item(554, 10).Value = 3000
$ws.Cells.Item(554, 11).Value = 500
$ws.Cells.Item(554, 12).Value = 500
$ws.Cells.Item(554, 13).Value = 500
$ws.Cells.Item(554, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(554, 15).Value = "Región del Maule"
$ws.Cells.Item(554, 16).Value = 100
$ws.Cells.Item(554, 17).Value = 5
$ws.Cells.Item(554, 18).Value = "Hortaliza"
